$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("H:H").Delete()
